$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 72: 2025-12-15
$ws.Cells.Item(72, 1).NumberFormat = "@"
$ws.Cells.Item(72, 1).Value = "2025-12-15"
$ws.Cells.Item(72, 1).ClearFormats()
$ws.Cells.Item(72, 2).Value = 0
$ws.Cells.Item(72, 3).Value = 31

# Row 73: 2025-12-16
$ws.Cells.Item(73, 1).NumberFormat = "@"
$ws.Cells.Item(73, 1).Value = "2025-12-16"
$ws.Cells.Item(73, 1).ClearFormats()
$ws.Cells.Item(73, 2).Value = 0
$ws.Cells.Item(73, 3).Value = 32
